$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.922.75'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.465.07'
$ws.Range('E3').Value = '  -2.69%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.25'
$ws.Range('E5').Value = '  -1.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.16'
$ws.Range('E6').Value = '  -1.73%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -2.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.463.55'
$ws.Range('E9').Value = '  -2.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.135'
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('E12').Value = '  -2.10%  '
$ws.Range('E13').Value = '  -4.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.911.78'
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.11'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.786.04'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000168'
$ws.Range('E17').Value = '  -3.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.396.71'
$ws.Range('E18').Value = '  -5.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.93'
$ws.Range('E19').Value = '  -6.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.36'
$ws.Range('E20').Value = '  -8.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '348.66'
$ws.Range('E21').Value = '  -4.66%  '
$ws.Range('E22').Value = '  -3.66%  '
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.36'
$ws.Range('E24').Value = '  -4.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.18'
$ws.Range('E25').Value = '  -7.33%  '
$ws.Range('E26').Value = '  -3.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.06'
$ws.Range('E27').Value = '  -7.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -37.56%  '
$ws.Range('E29').Value = '  -2.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0892'
$ws.Range('E30').Value = '  -5.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '506.22'
$ws.Range('E31').Value = '  -5.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.59'
$ws.Range('E32').Value = '  -7.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.76'
$ws.Range('E33').Value = '  -5.28%  '
$ws.Range('E34').Value = '  -5.38%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '158.29'
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('E37').Value = '  -11.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.66'
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.16'
$ws.Range('E39').Value = '  -5.67%  '
$ws.Range('E40').Value = '  -7.96%  '
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('E42').Value = '  -5.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.77'
$ws.Range('E43').Value = '  -5.52%  '
$ws.Range('E44').Value = '  -5.57%  '
$ws.Range('E45').Value = '  -4.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.81'
$ws.Range('E46').Value = '  -1.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '140.58'
$ws.Range('E47').Value = '  -4.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.42'
$ws.Range('E48').Value = '  -7.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.507'
$ws.Range('E49').Value = '  -7.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₆0249'
$ws.Range('E50').Value = '  -8.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0728'
$ws.Range('E51').Value = '  -2.03%  '
